# Apply the "Working on implementing assets" edit to the AssetList sheet.
#
# Summary of the change:
#  - Rows 13, 18, 19, 20: "Status of Implementation" flips from Incomplete
#    (red) to Complete (green) - the formula/dialog/music events are done.
#  - The single "EnemyAmbiance" row (covering both Vampire & Hellhound) is
#    split into two separate rows: VampireAmbiance and HellhoundAmbiance.
#    This pushes PlayerAmbiance, HealingTileAmbiance, Exit Sign Buzz and
#    Staff Fire Crackling down by one row each.
#  - HealingTileAmbiance and Exit Sign Buzz flip to Complete.
#  - A brand new row is appended for "Staff Fire Crackling" in its new
#    position (row 26), also marked Complete.
#  - Several Notes (column H) are updated to reflect current progress.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$COMPLETE_COLOR = 5296274   # RGB(0x92,0xD0,0x50) light green fill -> "Complete"
$INCOMPLETE_COLOR = 255     # RGB(0xFF,0x00,0x00) red fill -> "Incomplete"

function Set-Status($row, $col, $complete) {
    $cell = $ws.Cells.Item($row, $col)
    if ($complete) {
        $cell.Value = "Complete"
        $cell.Interior.Color = $COMPLETE_COLOR
    } else {
        $cell.Value = "Incomplete"
        $cell.Interior.Color = $INCOMPLETE_COLOR
    }
}

# --- Rows whose Status of Implementation (column G) became Complete ---
Set-Status 13 7 $true
Set-Status 18 7 $true
Set-Status 19 7 $true
Set-Status 20 7 $true

# --- Re-lay the Ambiance section (rows 21-26) ---

# Row 21: was EnemyAmbiance (Vampire+Hellhound) -> now VampireAmbiance alone
$ws.Cells.Item(21, 1).Value = "VampireAmbiance"
$ws.Cells.Item(21, 2).Value = "Ambiance for the enemies. Will include ambiance for each different type of enemy. "
$ws.Cells.Item(21, 3).Value = "Ambiance"
$ws.Cells.Item(21, 4).Value = "Vampire Ambiance"
Set-Status 21 5 $true
Set-Status 21 6 $true
Set-Status 21 7 $false
$ws.Cells.Item(21, 8).Value = "Having trouble because it doesn't recognize that there is a parameter"

# Row 22: new HellhoundAmbiance row (the other half of the old EnemyAmbiance row)
$ws.Cells.Item(22, 1).Value = "HellhoundAmbiance"
$ws.Cells.Item(22, 2).Value = "Ambiance for the enemies. Will include ambiance for each different type of enemy. "
$ws.Cells.Item(22, 3).Value = "Ambiance"
$ws.Cells.Item(22, 4).Value = "Hellhound Ambiance"
Set-Status 22 5 $true
Set-Status 22 6 $true
Set-Status 22 7 $false
$ws.Cells.Item(22, 8).Value = "Having trouble because it doesn't recognize that there is a parameter"

# Row 23: PlayerAmbiance (shifted down from row 22, note text changed)
$ws.Cells.Item(23, 1).Value = "PlayerAmbiance"
$ws.Cells.Item(23, 2).Value = "Ambiance for the player such as footsteps and breathing. "
$ws.Cells.Item(23, 3).Value = "Ambiance"
$ws.Cells.Item(23, 4).Value = "Player Ambiance"
Set-Status 23 5 $true
Set-Status 23 6 $true
Set-Status 23 7 $false
$ws.Cells.Item(23, 8).Value = "Having trouble because it doesn't recognize that there is a parameter"

# Row 24: HealingTileAmbiance (shifted down from row 23, now Complete)
$ws.Cells.Item(24, 1).Value = "HealingTileAmbiance"
$ws.Cells.Item(24, 2).Value = "Ambiance for the healing tiles. "
$ws.Cells.Item(24, 3).Value = "Ambiance"
$ws.Cells.Item(24, 4).Value = "Healing Tile Ambiance"
Set-Status 24 5 $true
Set-Status 24 6 $true
Set-Status 24 7 $true
$ws.Cells.Item(24, 8).Value = "Changed to 3D event. Currently always playing if there is a healing tile in the level, but I'm not sure why. "

# Row 25: Exit Sign Buzz (shifted down from row 24, now Complete, new note)
$ws.Cells.Item(25, 1).Value = "Exit Sign Buzz"
$ws.Cells.Item(25, 2).Value = "A sound for the buzz of the exit sign at the end of the game"
$ws.Cells.Item(25, 3).Value = "Ambiance"
$ws.Cells.Item(25, 4).Value = "Exit Sign Buzz"
Set-Status 25 5 $true
Set-Status 25 6 $true
Set-Status 25 7 $true
$ws.Cells.Item(25, 8).Value = "Is playing, but volume needs some attention in mixer"

# Row 26: Staff Fire Crackling - brand new row (shifted down from row 25, now Complete, new note)
$ws.Cells.Item(26, 1).Value = "Staff Fire Crackling"
$ws.Cells.Item(26, 2).Value = "A sound for the fire crackling on the player's staff in the menu and end screens"
$ws.Cells.Item(26, 3).Value = "Interface"
$ws.Cells.Item(26, 4).Value = "Fire Crackling"
Set-Status 26 5 $true
Set-Status 26 6 $true
Set-Status 26 7 $true
$ws.Cells.Item(26, 8).Value = "Needs volume attention with mixer"
